$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Regulation" dropdown value in B5 was re-entered in lowercase
# ("CSVLA" -> "csvla").
$ws.Range("B5").Value = "csvla"

# The underlying dropdown list (Data Validation) for B5 was also
# updated to use the lowercase options.
$dv = $ws.Range("B5").Validation
$dv.Formula1 = """csvla, cs-23, cs-25"""

# Selection/view was left on B5 after the edit.
$ws.Range("B5").Select()
